$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the TC_ID prefix used in this file
$ws.Name = "SCD0019"

# Update the TC_ID cell value (B2) from "DGS-321" to "SCD0019-004"
$ws.Range("B2").Value = "SCD0019-004"

# Let column B re-fit its width to the new (longer) TC_ID text
$ws.Columns.Item(2).AutoFit()

# Update the active selection to B3 as recorded in the saved view state
$ws.Range("B3").Select()
